# Apply the commit's changes to the presentation:
#  1. Update the title text on slide 1 (shorten it).
#  2. Delete slide 5 ("Przygotowanie danych"), which also shifts every
#     later slide up by one position (handled automatically by Delete()).

$p = $ppt.ActivePresentation

# --- 1. Slide 1: shorten the title -----------------------------------
$titleSlide = $p.Slides.Item(1)
$titleShape = $titleSlide.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleRange.Text = "Przewidywanie opadów"

# Re-assign per-word so the new title keeps the same run breaks
# ("Przewidywanie" / " " / "opadów") as the authored slide.
$titleRange.Characters(1, 13).Text = "Przewidywanie"
$titleRange.Characters(14, 1).Text = " "
$titleRange.Characters(15, 6).Text = "opadów"

# --- 2. Delete slide 5 ("Przygotowanie danych") -----------------------
$obsoleteSlide = $p.Slides.Item(5)
$obsoleteSlide.Delete()
